# Saldo.xlsx update
#
# Source data (Export sheet) is a flat "Conta / Nome / Saldo" export.
# This change:
#   1) Adds a new account row for LUIZ (004458563, 300000) at the top of
#      the data (right after the header row).
#   2) Moves the CLUBE account (005696595) up near the top of the list and
#      updates its balance from 87.11 to 17963.71 (i.e. removes the old
#      low-balance row further down and inserts a new row, with the new
#      balance, right before the MURYLO row).
#   3) Removes the MARIANA account (005000460) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the old CLUBE row (005696595 / CLUBE / 87.11), further down
#         the list, between VICTOR (003115072) and JOAO (004212132). ---
$ws.Rows.Item(37).Delete()

# --- 2) Remove the MARIANA row (005000460 / MARIANA / 2000). ---
$ws.Rows.Item(5).Delete()

# --- 3) Insert the new LUIZ row right after the header, before GUSTAVO. ---
#         (NumberFormat "@" forces the leading-zero account number to be
#         kept as text instead of Excel auto-converting it to a number;
#         ClearFormats afterwards drops the now-unneeded format override so
#         the cell ends up plain/unstyled, like the rest of the sheet, while
#         remaining text.)
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004458563"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "LUIZ"
$ws.Cells.Item(2, 3).Value = 300000

# --- 4) Insert the new CLUBE row (with updated balance) right before the
#         MURYLO row (004397124). ---
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "005696595"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "CLUBE"
$ws.Cells.Item(5, 3).Value = 17963.71
